$d = $word.ActiveDocument

$d.Content.Find.Execute("265×4=1060", $true, $false, $false, $false, $false, $true, 1, $false, "982×6=5892", 2)
$d.Content.Find.Execute("114×7=798", $true, $false, $false, $false, $false, $true, 1, $false, "348×7=2436", 2)
$d.Content.Find.Execute("245×9=2205", $true, $false, $false, $false, $false, $true, 1, $false, "333×7=2331", 2)
$d.Content.Find.Execute("643×4=2572", $true, $false, $false, $false, $false, $true, 1, $false, "485×6=2910", 2)
$d.Content.Find.Execute("495×4=1980", $true, $false, $false, $false, $false, $true, 1, $false, "236×9=2124", 2)
$d.Content.Find.Execute("307×7=2149", $true, $false, $false, $false, $false, $true, 1, $false, "187×3=561", 2)
$d.Content.Find.Execute("979×9=8811", $true, $false, $false, $false, $false, $true, 1, $false, "402×8=3216", 2)
$d.Content.Find.Execute("648×4=2592", $true, $false, $false, $false, $false, $true, 1, $false, "341×7=2387", 2)
$d.Content.Find.Execute("734×5=3670", $true, $false, $false, $false, $false, $true, 1, $false, "117×9=1053", 2)
$d.Content.Find.Execute("826×2=1652", $true, $false, $false, $false, $false, $true, 1, $false, "195×2=390", 2)
$d.Content.Find.Execute("836×4=3344", $true, $false, $false, $false, $false, $true, 1, $false, "482×3=1446", 2)
$d.Content.Find.Execute("261×4=1044", $true, $false, $false, $false, $false, $true, 1, $false, "338×6=2028", 2)
$d.Content.Find.Execute("303×6=1818", $true, $false, $false, $false, $false, $true, 1, $false, "731×9=6579", 2)
$d.Content.Find.Execute("694×2=1388", $true, $false, $false, $false, $false, $true, 1, $false, "496×7=3472", 2)
$d.Content.Find.Execute("139×8=1112", $true, $false, $false, $false, $false, $true, 1, $false, "744×9=6696", 2)
$d.Content.Find.Execute("434×4=1736", $true, $false, $false, $false, $false, $true, 1, $false, "413×8=3304", 2)
$d.Content.Find.Execute("580×5=2900", $true, $false, $false, $false, $false, $true, 1, $false, "895×2=1790", 2)
$d.Content.Find.Execute("607×5=3035", $true, $false, $false, $false, $false, $true, 1, $false, "667×7=4669", 2)
$d.Content.Find.Execute("424×3=1272", $true, $false, $false, $false, $false, $true, 1, $false, "413×5=2065", 2)
$d.Content.Find.Execute("700×9=6300", $true, $false, $false, $false, $false, $true, 1, $false, "320×7=2240", 2)
$d.Content.Find.Execute("160×5=800", $true, $false, $false, $false, $false, $true, 1, $false, "608×7=4256", 2)
$d.Content.Find.Execute("638×4=2552", $true, $false, $false, $false, $false, $true, 1, $false, "105×3=315", 2)
$d.Content.Find.Execute("272×7=1904", $true, $false, $false, $false, $false, $true, 1, $false, "844×5=4220", 2)
$d.Content.Find.Execute("316×8=2528", $true, $false, $false, $false, $false, $true, 1, $false, "219×4=876", 2)
$d.Content.Find.Execute("926×8=7408", $true, $false, $false, $false, $false, $true, 1, $false, "317×6=1902", 2)
